$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

function Set-PlainValue($cellRef, $value) {
    $ws.Range($cellRef).Value = $value
}

Set-PlainValue 'D2' '62.370.83'
Set-PlainValue 'E2' '  -1.91%  '
Set-PlainValue 'D3' '3.011.12'
Set-PlainValue 'E3' '  -1.88%  '
Set-TextValue 'D4' '1.00'
Set-PlainValue 'E4' '  +0.08%  '
Set-TextValue 'D5' '583.04'
Set-PlainValue 'E5' '  -0.60%  '
Set-TextValue 'D6' '147.14'
Set-PlainValue 'E6' '  -4.51%  '
Set-PlainValue 'E7' '  +0.03%  '
Set-PlainValue 'B8' 'LidoStakedEther'
Set-PlainValue 'C8' 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
Set-PlainValue 'D8' '3.014.82'
Set-PlainValue 'E8' '  -1.70%  '
Set-PlainValue 'B9' 'XRP'
Set-PlainValue 'C9' 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextValue 'D9' '0.520'
Set-PlainValue 'E9' '  -3.11%  '
Set-PlainValue 'E10' '  -4.41%  '
Set-TextValue 'D11' '5.69'
Set-PlainValue 'E11' '  -2.06%  '
Set-TextValue 'D12' '0.441'
Set-PlainValue 'E12' '  -1.44%  '
Set-PlainValue 'E13' '  -3.01%  '
Set-TextValue 'D14' '34.73'
Set-PlainValue 'E14' '  -5.06%  '
Set-TextValue 'D15' '0.122'
Set-PlainValue 'E15' '  +2.23%  '
Set-PlainValue 'D16' '3.514.71'
Set-PlainValue 'E16' '  -1.64%  '
Set-TextValue 'D17' '7.05'
Set-PlainValue 'E17' '  -0.92%  '
Set-PlainValue 'D18' '62.433.09'
Set-PlainValue 'E18' '  -1.70%  '
Set-PlainValue 'D19' '3.014.92'
Set-PlainValue 'E19' '  -1.77%  '
Set-TextValue 'D20' '459.98'
Set-PlainValue 'E20' '  -1.87%  '
Set-PlainValue 'E21' '  -2.24%  '
Set-TextValue 'D22' '0.684'
Set-TextValue 'D23' '7.33'
Set-PlainValue 'E23' '  -1.86%  '
Set-PlainValue 'E24' '  -6.27%  '
Set-TextValue 'D25' '79.99'
Set-PlainValue 'E25' '  -0.25%  '
Set-PlainValue 'E26' '  -3.53%  '
Set-PlainValue 'B27' 'RenderToken'
Set-PlainValue 'C27' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D27' '10.09'
Set-PlainValue 'E27' '  -2.51%  '
Set-PlainValue 'B28' 'Dai'
Set-PlainValue 'C28' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D28' '0.999'
Set-PlainValue 'E28' '  -0.08%  '
Set-PlainValue 'E29' '  -0.04%  '
Set-PlainValue 'E30' '  -1.09%  '
Set-TextValue 'D31' '7.15'
Set-PlainValue 'E31' '  -2.52%  '
Set-PlainValue 'E32' '  -1.52%  '
Set-TextValue 'D33' '27.09'
Set-PlainValue 'E33' '  +0.33%  '
Set-PlainValue 'E34' '  -4.44%  '
Set-PlainValue 'E35' '  -0.76%  '
Set-PlainValue 'D36' '0.0₃0789'
Set-PlainValue 'E36' '  -3.52%  '
Set-TextValue 'D37' '5.74'
Set-PlainValue 'E37' '  -3.46%  '
Set-PlainValue 'E38' '  -3.84%  '
Set-TextValue 'D39' '50.39'
Set-PlainValue 'E39' '  -0.15%  '
Set-TextValue 'D40' '9.00'
Set-PlainValue 'E40' '  -1.05%  '
Set-TextValue 'D41' '2.89'
Set-PlainValue 'E41' '  -10.52%  '
Set-TextValue 'D42' '415.84'
Set-PlainValue 'E42' '  -4.77%  '
Set-PlainValue 'E43' '  +0.96%  '
Set-PlainValue 'E44' '  -4.06%  '
Set-PlainValue 'B45' 'VeChain'
Set-PlainValue 'C45' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D45' '0.0354'
Set-PlainValue 'E45' '  -1.12%  '
Set-PlainValue 'B46' 'Maker'
Set-PlainValue 'C46' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-PlainValue 'D46' '2.774.53'
Set-PlainValue 'E46' '  -0.76%  '
Set-TextValue 'D47' '38.06'
Set-PlainValue 'E47' '  -5.14%  '
Set-TextValue 'D48' '128.59'
Set-PlainValue 'E48' '  -1.32%  '
Set-TextValue 'D49' '1.00'
Set-PlainValue 'E49' '  +0.00%  '
Set-PlainValue 'E50' '  -0.74%  '
Set-TextValue 'D51' '23.80'
Set-PlainValue 'E51' '  -4.14%  '
